# Dragon Discounts: Content - Added segmentation by player source.
# Adds a new "[playerSources]" column to the end of Table1 on Sheet1,
# filling its data rows with the placeholder "-" value used by the
# other not-yet-configured segmentation columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add a new (40th) column to the table; Excel appends it right after
# the current last column ([clusterId], col AN) as new col AO.
$newCol = $lo.ListColumns.Add()
$hdrCell = $newCol.Range.Item(1)

# Match formatting: the new last header cell picks up the "regular"
# repeating header style (same as the header one step to its left),
# while the data cells copy the style banding already used by the
# previous last column (AN) so the row-striping stays consistent.
$ws.Range("AM2").Copy() | Out-Null
$hdrCell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("AN3:AN10").Copy() | Out-Null
$ws.Range("AO3:AO10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Header text
$hdrCell.Value = "[playerSources]"

# Data rows: default placeholder "-" (same as sibling segmentation columns)
for ($r = 3; $r -le 10; $r++) {
    $ws.Range("AO$r").Value = "-"
}

# Cosmetic: widen the new column and move the view/selection the way the
# authored workbook ended up (not data-bearing, best effort only).
$ws.Range("AO1").ColumnWidth = 27.75
$ws.Range("AO22").Select() | Out-Null

Write-Host "Added [playerSources] column to Table1 (now $($lo.Range.Address()))"
